# Update cryptocurrency price and volume figures to refreshed values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.252.08'
$ws.Range('E2').Value = '  +4.07%  '
$ws.Range('D3').Value = '1.727.43'
$ws.Range('E3').Value = '  +2.80%  '
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '219.38'
$ws.Range('E5').Value = '  +1.91%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.521'
$ws.Range('E6').Value = '  +0.63%  '
$ws.Range('E7').Value = '  -0.26%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '24.35'
$ws.Range('E8').Value = '  +13.98%  '
$ws.Range('E9').Value = '  +3.30%  '
$ws.Range('E10').Value = '  +1.81%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0900'
$ws.Range('E11').Value = '  +1.33%  '
$ws.Range('D12').Value = '1.969.45'
$ws.Range('E12').Value = '  +2.73%  '
$ws.Range('D13').Value = '1.727.10'
$ws.Range('E13').Value = '  +2.83%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.27'
$ws.Range('E14').Value = '  +3.09%  '
$ws.Range('E15').Value = '  +4.61%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '67.65'
$ws.Range('E16').Value = '  +2.04%  '
$ws.Range('D17').Value = '28.188.57'
$ws.Range('E17').Value = '  +3.84%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '244.27'
$ws.Range('E18').Value = '  +2.05%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '8.04'
$ws.Range('E19').Value = '  -0.48%  '
$ws.Range('D20').Value = '0.0₃0756'
$ws.Range('E20').Value = '  +1.78%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.999'
$ws.Range('E21').Value = '  -0.21%  '
$ws.Range('E22').Value = '  +2.73%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.69'
$ws.Range('E23').Value = '  +2.49%  '
$ws.Range('E24').Value = '  -0.58%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '149.43'
$ws.Range('E25').Value = '  +1.59%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.52'
$ws.Range('E26').Value = '  +3.87%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '16.74'
$ws.Range('E27').Value = '  +2.28%  '
$ws.Range('E28').Value = '  +0.89%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.27%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0512'
$ws.Range('E30').Value = '  +2.55%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.20'
$ws.Range('E31').Value = '  +2.13%  '
$ws.Range('E32').Value = '  +2.30%  '
$ws.Range('D33').Value = '1.505.08'
$ws.Range('E33').Value = '  -3.93%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.28'
$ws.Range('E34').Value = '  +1.92%  '
$ws.Range('E35').Value = '  -1.24%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.965'
$ws.Range('E36').Value = '  +3.52%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.610'
$ws.Range('E37').Value = '  +1.49%  '
$ws.Range('E38').Value = '  +0.30%  '
$ws.Range('E39').Value = '  +0.84%  '
$ws.Range('E40').Value = '  +1.26%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '70.96'
$ws.Range('E41').Value = '  +2.57%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.79'
$ws.Range('E42').Value = '  +4.09%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.999'
$ws.Range('E43').Value = '  -0.20%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.30'
$ws.Range('E44').Value = '  +1.83%  '
$ws.Range('D45').Value = '1.874.92'
$ws.Range('E45').Value = '  +2.52%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.807'
$ws.Range('E46').Value = '  +2.63%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.78'
$ws.Range('E47').Value = '  +12.12%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '90.91'
$ws.Range('E48').Value = '  +0.26%  '
$ws.Range('D49').Value = '0.0₆0113'
$ws.Range('E49').Value = '  +6.39%  '
$ws.Range('E50').Value = '  +1.30%  '
$ws.Range('E51').Value = '  +0.80%  '
